$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Item line: quantity & unit price change ---
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 10.95

# --- Relabel Sub Total row ---
$ws.Range("A18").Value = "Sub Total (VAT EXCLUSIVE)"

# --- Rework the totals formulas ---
$ws.Range("I18").Formula = "=(I15+I20)/1.2"
$ws.Range("I20").FormulaArray = "=IFS(G15=1, -I15*0, G15=2, -I15*0.1, G15=3,-I15*0.2)"
$ws.Range("I21").Formula = "=(SUM(I19:I20)+I15)*0.1667"
$ws.Range("I22").Formula = "=SUM(I21,I18)"

# --- Selection moves to the Sub Total row ---
$ws.Range("A18:H18").Select()
